$d = $word.ActiveDocument

# --- 1. Insert "Release: v1.0.0" paragraph right after the "Pull Request Template" heading (bold) ---
$pTitle = $d.Paragraphs.Item(1)
$pTitle.Range.InsertParagraphAfter()
$pRelease = $d.Paragraphs.Item(2)
$pRelease.Style = "FirstParagraph"
$rRelease = $pRelease.Range
$rRelease.Collapse(1)
$rRelease.InsertAfter("Release: v1.0.0")
$rRelease.Font.Bold = $true
$rRelease.Font.BoldBi = $true

# --- 2. Insert "Summary of the change and which issue is fixed." paragraph right after "Description" ---
$pDesc = $d.Paragraphs.Item(3)
$pDesc.Range.InsertParagraphAfter()
$pSummary = $d.Paragraphs.Item(4)
$pSummary.Style = "FirstParagraph"
$rSummary = $pSummary.Range
$rSummary.Collapse(1)
$rSummary.InsertAfter("Summary of the change and which issue is fixed.")

# --- 3. Checklist text replacements (shorten wording) ---
[void]$d.Content.Find.Execute(
    "My code follows the style guidelines of this project", $true, $false, $false, $false, $false,
    $true, 1, $false, "My code follows the style guidelines", 2)

[void]$d.Content.Find.Execute(
    "I have performed a self-review of my code", $true, $false, $false, $false, $false,
    $true, 1, $false, "I have performed a self-review", 2)

[void]$d.Content.Find.Execute(
    "I have commented my code, particularly in hard-to-understand areas", $true, $false, $false, $false, $false,
    $true, 1, $false, "I have commented my code", 2)

[void]$d.Content.Find.Execute(
    "I have made corresponding changes to the documentation", $true, $false, $false, $false, $false,
    $true, 1, $false, "I have made corresponding documentation changes", 2)

[void]$d.Content.Find.Execute(
    "I have added tests that prove my fix is effective or that my feature works", $true, $false, $false, $false, $false,
    $true, 1, $false, "I have added tests or proven my feature works", 2)

[void]$d.Content.Find.Execute(
    "New and existing unit tests pass locally with my changes", $true, $false, $false, $false, $false,
    $true, 1, $false, "All tests pass locally", 2)

# --- 4. Insert "Link to related issues or PRs" paragraph right after "Related Issues" (last paragraph) ---
$relatedIndex = $d.Paragraphs.Count
$pRelated = $d.Paragraphs.Item($relatedIndex)
$pRelated.Range.InsertParagraphAfter()
$pLink = $d.Paragraphs.Item($relatedIndex + 1)
$pLink.Style = "FirstParagraph"
$rLink = $pLink.Range
$rLink.Collapse(1)
$rLink.InsertAfter("Link to related issues or PRs")
